$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the email address in A3 (retry mechanism uses a new test account)
$ws.Range("A3").Value = "nasrath1298@binafex.com"

# Reuse the existing "Test@123" password for the new row (B3), replacing "NarrowBabe"
$ws.Range("B3").Value = "Test@123"

# Add a hyperlink for the retried password cell, same target as B2's hyperlink
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:Test@123")

# Adding the hyperlink applies the built-in "Hyperlink" cell style (underline/blue font);
# the target workbook keeps B3 on the default style, so reset it back to Normal.
$ws.Range("B3").Style = "Normal"

# Resize column A to a fixed width (no longer auto best-fit)
$ws.Columns.Item(1).ColumnWidth = 23.166666666666668

# Update the active cell / selection on the sheet
$ws.Range("D8").Select()
